$wb = $excel.ActiveWorkbook

# Sheet 1: ROW35-FE-LIFTER - append row 49
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(49,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(49,1).Value2 = 45748.84486297454
$ws1.Cells.Item(49,2).Value2 = "0x01,0x90"
$ws1.Cells.Item(49,3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item(49,4).Value2 = "0x01,0x76"
$ws1.Cells.Item(49,5).Value2 = "0xd"
$ws1.Cells.Item(49,6).Value2 = 400
$ws1.Cells.Item(49,7).Value2 = [double]"5.68631262647114e+23"
$ws1.Cells.Item(49,8).Value2 = 374
$ws1.Cells.Item(49,9).Value2 = 13

# Sheet 2: ROW35-MID-LIFTER - append row 49
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(49,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(49,1).Value2 = 45748.69673190972
$ws2.Cells.Item(49,2).Value2 = "0x01,0x90"
$ws2.Cells.Item(49,3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(49,4).Value2 = "0x01,0x72"
$ws2.Cells.Item(49,5).Value2 = "0xe"
$ws2.Cells.Item(49,6).Value2 = 400
$ws2.Cells.Item(49,7).Value2 = [double]"5.68631262647114e+23"
$ws2.Cells.Item(49,8).Value2 = 370
$ws2.Cells.Item(49,9).Value2 = 14

# Sheet 3: ROW02-FE-LIFTER - append row 49
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(49,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(49,1).Value2 = 45748.83747498842
$ws3.Cells.Item(49,2).Value2 = "0x01,0x90"
$ws3.Cells.Item(49,3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item(49,4).Value2 = "0x01,0x76"
$ws3.Cells.Item(49,5).Value2 = "0x3"
$ws3.Cells.Item(49,6).Value2 = 400
$ws3.Cells.Item(49,7).Value2 = [double]"5.68631262647114e+23"
$ws3.Cells.Item(49,8).Value2 = 374
$ws3.Cells.Item(49,9).Value2 = 3

# Sheet 4: ROW02-MID-LIFTER - append row 49
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(49,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item(49,1).Value2 = 45748.89468563657
$ws4.Cells.Item(49,2).Value2 = "0x01,0x90"
$ws4.Cells.Item(49,3).Value2 = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Cells.Item(49,4).Value2 = "0x01,0x72"
$ws4.Cells.Item(49,5).Value2 = "0x3"
$ws4.Cells.Item(49,6).Value2 = 400
$ws4.Cells.Item(49,7).Value2 = [double]"9.85046333984776e+23"
$ws4.Cells.Item(49,8).Value2 = 370
$ws4.Cells.Item(49,9).Value2 = 3
